$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, copying the formatting (style) of the
# existing H1 header cell, then overwrite the copied text with the new
# header labels.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "IF"

# Populate the new I/J data columns for rows 2-12 (no special style,
# matching the other plain numeric data cells).
$values = @{
  2  = @(1, 4)
  3  = @(4, 6)
  4  = @(8, 9)
  5  = @(4, 5)
  6  = @(6, 6)
  7  = @(4, 5)
  8  = @(8, 8)
  9  = @(9, 9)
  10 = @(7, 7)
  11 = @(6, 6)
  12 = @(5, 5)
}

foreach ($row in $values.Keys) {
  $pair = $values[$row]
  $ws.Cells.Item($row, 9).Value = $pair[0]
  $ws.Cells.Item($row, 10).Value = $pair[1]
}
